$wb = $excel.ActiveWorkbook

# --- Sheet "MERtoPPP": just a view/selection change (D2:D4 -> K2) ---
$mer = $wb.Worksheets.Item("MERtoPPP")
$mer.Range("K2").Select() | Out-Null

# --- Sheet "config": add a commodity column, reordering existing columns,
#     and drop the now-unused "year" column (which only had values in rows 2-4) ---
$cfg = $wb.Worksheets.Item("config")

# Insert a new blank column before the existing "level" column (C), shifting
# level -> D and commodity -> E.
$cfg.Range("C1").EntireColumn.Insert() | Out-Null

# Populate the freshly inserted column C with the "commodity" data (header + value).
$cfg.Range("C1").Value = "commodity"
$cfg.Range("C2").Value = "light"

# Clear out the old "commodity" column (now E) -- its data has been moved to C.
$cfg.Range("E1:E4").ClearContents() | Out-Null

# Clear out the "year" column (now F) -- it is no longer part of the config sheet.
$cfg.Range("F1:F4").ClearContents() | Out-Null

# Match the column width Excel computed when auto-fitting the new "commodity" column.
$cfg.Columns("C").ColumnWidth = 10.166666666666666

# Update the selected cell to match the saved view state.
$cfg.Range("A3").Select() | Out-Null
